$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 6: new test case "Correct product search" ----
$ws.Range("A6").Value = 'Correct product search'
$ws.Range("B6").Value = 'Poster V1'
$ws.Range("C6").Value = '-'
$ws.Range("D6").Value = '1. Open the website https://test.testowanie-oprogramowania.pl
2. Press the Shop tab
3. Insert into "Search products…" field the input value
4.Press the search button to the left of the field
'
$ws.Range("E6").Value = 'Redirection to the product page'
$ws.Range("F6").Value = 'POSITIVE'

# ---- Row 7: new test case "Invalid coupon code" ----
$ws.Range("A7").Value = 'Invalid coupon code'
$ws.Range("B7").Value = 'xyz321u5s'
$ws.Range("C7").Value = 'Product is in the shopping cart'
$ws.Range("D7").Value = '1. Open the website https://test.testowanie-oprogramowania.pl
2. Press the Shopping cart tab
3. Press the "Zobacz koszyk" button
4.Insert into "Kod kuponu" field input value
5. Press the "Wykorzystaj kupon" button
'
$ws.Range("E7").Value = 'Displayed message: "Kupon "xyz321u5s" nie istnieje!"'
$ws.Range("F7").Value = 'POSITIVE'

# ---- Styles: center the "-" cell, copy POSITIVE formatting, set wrap/top for D column ----
$ws.Range("C6").HorizontalAlignment = -4108

$ws.Range("F2").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F7").PasteSpecial(-4122)

$ws.Range("A7").Copy()
$ws.Range("E7").PasteSpecial(-4122)

# D2:D7 get the new "wrap + vertical top" style
$ws.Range("D2").WrapText = $true
$ws.Range("D2").VerticalAlignment = -4160
$ws.Range("D2").Copy()
$ws.Range("D3:D7").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- Row heights ----
$ws.Rows.Item(3).RowHeight = 122.4
$ws.Rows.Item(4).RowHeight = 114
$ws.Rows.Item(5).RowHeight = 113.4
$ws.Rows.Item(6).RowHeight = 124.2
$ws.Rows.Item(7).RowHeight = 153.6

# ---- Column widths ----
$ws.Columns.Item(4).ColumnWidth = 33.75
$ws.Columns.Item(5).ColumnWidth = 48.25

# ---- Selection matches the author's saved cursor position ----
$ws.Range("D10").Select()
